# Update the "想去人数" (F column) counts on the three worksheets that carry
# this data: 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4).
# 本地生活 (sheet3) has no data rows and is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2793
$ws1.Range("F6").Value  = 2445
$ws1.Range("F7").Value  = 15
$ws1.Range("F11").Value = 37
$ws1.Range("F13").Value = 6973
$ws1.Range("F14").Value = 272
$ws1.Range("F16").Value = 213
$ws1.Range("F18").Value = 469
$ws1.Range("F19").Value = 8089
$ws1.Range("F22").Value = 258
$ws1.Range("F27").Value = 65
$ws1.Range("F28").Value = 32
$ws1.Range("F33").Value = 2591
$ws1.Range("F34").Value = 39
$ws1.Range("F35").Value = 82
$ws1.Range("F38").Value = 87
$ws1.Range("F39").Value = 636
$ws1.Range("F40").Value = 3647
$ws1.Range("F41").Value = 164
$ws1.Range("F42").Value = 1171
$ws1.Range("F43").Value = 144

# --- Sheet "演出" ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 2
$ws2.Range("F17").Value = 32

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2793
$ws4.Range("F8").Value  = 2445
$ws4.Range("F10").Value = 15
$ws4.Range("F14").Value = 37
$ws4.Range("F18").Value = 6973
$ws4.Range("F19").Value = 272
$ws4.Range("F21").Value = 213
$ws4.Range("F23").Value = 469
$ws4.Range("F24").Value = 8089
$ws4.Range("F27").Value = 258
$ws4.Range("F32").Value = 65
$ws4.Range("F38").Value = 2591
$ws4.Range("F39").Value = 39
$ws4.Range("F40").Value = 82
$ws4.Range("F43").Value = 87
$ws4.Range("F44").Value = 636
$ws4.Range("F46").Value = 3647
$ws4.Range("F47").Value = 164
$ws4.Range("F49").Value = 1171
$ws4.Range("F50").Value = 144

$wb.Save()
